# DPLKINV129-001 / DPLKINV129-002 - Fixed Income Approve Dealing Ticket
# Update the Ticket IDs referenced in the two test-case sheets (new run,
# new dealing tickets DTOBL202300026 / DTOBL202300027), and move the
# active sheet/selection/scroll position like the author left the
# workbook after editing.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # DPLKINV129-001
$ws2 = $wb.Worksheets.Item(2)   # DPLKINV129-002

# --- Sheet 1: DPLKINV129-001 (Verifikasi disetujui) ---------------------
[void]$ws1.Activate()

$ws1.Range("N2").Value = "DTOBL202300026"
$ws1.Range("F2").Value = "Username : 33372;`nPassword : bni1234;`nRole : 18 - Pimpinan Kelompok Investasi;`nTicket ID : DTOBL202300026;`nStatus Verifikasi : 1 : Setuju;`nKeterangan Verifikasi : INV.FIX.BEJ.011 DISETUJUI"

$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
[void]$ws1.Range("G2").Select()

# --- Sheet 2: DPLKINV129-002 (Kembalikan ke Data Entry) ------------------
[void]$ws2.Activate()

$ws2.Range("N2").Value = "DTOBL202300027"
$ws2.Range("F2").Value = "Username : 33372;`nPassword : bni1234;`nRole : 18 - Pimpinan Kelompok Investasi;`nTicket ID : DTOBL202300027;`nStatus Verifikasi : 0 : Kembalikan ke Data Entry;`nKeterangan Verifikasi : INV.FIX.BEJ.011 dikembalikan ke Data Entry"

$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
[void]$ws2.Range("G2").Select()

# DPLKINV129-002 stays the active (last-activated / tabSelected) sheet,
# matching the workbook's new activeTab="1".
